# Update "API Access Control - Concern" workbook:
#  - rename "Sheet1" to "Data" (the _FilterDatabase defined name and the
#    AutoFilter range follow the rename automatically)
#  - turn the Legend sheet's A1:B6 lookup range into a header row + a real
#    Excel Table ("Table1") with generic header labels "Column1"/"Column2"
#  - re-apply Data!H1's header style to Data!I1:J1 so the two trailing
#    question columns share the same formatting as the rest of the header row

$wb = $excel.ActiveWorkbook

$wsData   = $wb.Worksheets.Item("Sheet1")
$wsLegend = $wb.Worksheets.Item("Legend")

# 1) Rename the data sheet.
$wsData.Name = "Data"

# 2) Insert a header row above the existing Legend lookup rows and label the
#    two columns, then convert the range into a Table (ListObject) so the
#    sheet carries a proper header + table, matching the authored workbook.
$wsLegend.Range("A1").EntireRow.Insert() | Out-Null
$wsLegend.Range("A1").Value = "Column1"
$wsLegend.Range("B1").Value = "Column2"

$legendTable = $wsLegend.ListObjects.Add(1, $wsLegend.Range("A1:B7"), 0, 1)
$legendTable.Name = "Table1"

# Leave the whole table selected, as Excel does right after "Format as Table".
$wsLegend.Range("A1:B7").Select() | Out-Null

# 3) Line I1/J1 back up with H1's header style so all three trailing
#    question-header cells share one consistent cell format.
$wsData.Range("H1").Copy() | Out-Null
$wsData.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Restore the original active sheet/selection context.
$wsData.Activate() | Out-Null
$wsData.Range("H23").Select() | Out-Null
